$d = $word.ActiveDocument

$replacements = @(
    @{old = "2023-12-06 Wednesday"; new = "2023-12-07 Thursday"},
    @{old = "81×47="; new = "96×96="},
    @{old = "46×28="; new = "32×71="},
    @{old = "47×17="; new = "94×55="},
    @{old = "18×20="; new = "18×94="},
    @{old = "56×81="; new = "57×49="},
    @{old = "45×67="; new = "96×80="},
    @{old = "73×34="; new = "39×30="},
    @{old = "52×17="; new = "84×99="},
    @{old = "29×55="; new = "86×66="},
    @{old = "85×27="; new = "17×83="},
    @{old = "89×82="; new = "69×86="},
    @{old = "84×80="; new = "53×48="},
    @{old = "28×61="; new = "20×78="},
    @{old = "80×65="; new = "48×64="},
    @{old = "90×55="; new = "32×47="},
    @{old = "55×79="; new = "87×90="},
    @{old = "14×17="; new = "36×93="},
    @{old = "32×42="; new = "77×94="},
    @{old = "11×86="; new = "15×30="},
    @{old = "12×92="; new = "13×15="},
    @{old = "87×43="; new = "56×64="},
    @{old = "23×28="; new = "94×16="},
    @{old = "76×87="; new = "54×34="},
    @{old = "64×79="; new = "76×58="},
    @{old = "69×73="; new = "40×79="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
